# blocktypes.xlsx update: "L41, and loads of new assets"
#
# Adds a second "blocks" table (columns C/D, mirroring the A/B "Block type" /
# "block letter" table) with the new Carpet/cratE/canDelabrum/cHest rows, and
# adds an "enemies" / "blocks" letter-reference table in columns H/I listing
# every letter A-Z split between the letters used by enemies and the letters
# used by blocks (the newly-added block letters C, D, E and H are shown in
# bold in the blocks column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate new cells --------------------------------------------------
# NOTE: the order in which brand-new text values are first written matters,
# since it determines the order new entries are appended to the shared
# string table. The order below reproduces the target shared-string layout.

# New block-letter column (D) - first occurrences of each new single-letter
# string, interleaved with a couple of "blocks" column (I) letters so that
# the overall shared-string order comes out as C, D, E, F, H, I, J, L, N, Q,
# V, X, Y, Z, enemies, blocks, cratE, canDelabrum, cHest, Carpet.
$ws.Range("D4").Value = "C"
$ws.Range("D6").Value = "D"
$ws.Range("D5").Value = "E"
$ws.Range("I9").Value = "F"
$ws.Range("D7").Value = "H"
$ws.Range("I12").Value = "I"
$ws.Range("I13").Value = "J"
$ws.Range("I15").Value = "L"
$ws.Range("I17").Value = "N"
$ws.Range("I20").Value = "Q"
$ws.Range("I25").Value = "V"
$ws.Range("I27").Value = "X"
$ws.Range("I28").Value = "Y"
$ws.Range("I29").Value = "Z"

# Header row for the new enemies / blocks reference table.
$ws.Range("H3").Value = "enemies"
$ws.Range("H3").Font.Bold = $true
$ws.Range("I3").Value = "blocks"
$ws.Range("I3").Font.Bold = $true

# New "block type" rows (column C), mirroring columns A/B.
$ws.Range("C5").Value = "cratE"
$ws.Range("C6").Value = "canDelabrum"
$ws.Range("C7").Value = "cHest"
$ws.Range("C4").Value = "Carpet"

# Header row for the new block-type table (columns C/D), same text as A3/B3.
$ws.Range("C3").Value = "Block type"
$ws.Range("C3").Font.Bold = $true
$ws.Range("D3").Value = "block letter"
$ws.Range("D3").Font.Bold = $true

# Remaining "enemies" column (H) letters (reuse of existing single-letter
# shared strings already present in the workbook).
$ws.Range("H4").Value = "A"
$ws.Range("H5").Value = "B"
$ws.Range("H10").Value = "G"
$ws.Range("H14").Value = "K"
$ws.Range("H16").Value = "M"
$ws.Range("H18").Value = "O"
$ws.Range("H19").Value = "P"
$ws.Range("H21").Value = "R"
$ws.Range("H22").Value = "S"
$ws.Range("H23").Value = "T"
$ws.Range("H24").Value = "U"
$ws.Range("H26").Value = "W"

# Remaining "blocks" column (I) letters that are also the block letters used
# above (C, D, E, H) are shown in bold.
$ws.Range("I6").Value = "C"
$ws.Range("I6").Font.Bold = $true
$ws.Range("I7").Value = "D"
$ws.Range("I7").Font.Bold = $true
$ws.Range("I8").Value = "E"
$ws.Range("I8").Font.Bold = $true
$ws.Range("I11").Value = "H"
$ws.Range("I11").Font.Bold = $true

# --- Column widths ---------------------------------------------------------
# New custom widths for the two added columns (closest values reachable
# through the ColumnWidth property, which snaps to whole-pixel increments).
$ws.Columns.Item(3).ColumnWidth = 25.8
$ws.Columns.Item(4).ColumnWidth = 20.6

# --- View state --------------------------------------------------------
# Scroll so row 4 is at the top and select C20, matching the saved view.
$ws.Range("A4").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C20").Select() | Out-Null
